$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2 through 57 is updated from serial date
# 45188 (2023-09-19) to 45189 (2023-09-20).
$ws.Range("C2:C57").Value = 45189
